$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing row (42) onto the new row (43)
# so the new row inherits the same cell styles (fonts/borders/alignment).
$ws.Range("A42:H42").Copy()
$ws.Range("A43:H43").PasteSpecial(-4122)   # xlPasteFormats

# Populate the new row's values (mirrors most columns from row 42, with new
# Status Date / Fail KPI / Test Case / Remarks content).
$ws.Range("A43").Value = "KK"
$ws.Range("B43").Value = "GWPRA1_DAV_P47"
$ws.Range("C43").Value = "L2100"
$ws.Range("D43").Value = "23-Dec-2025 10:04 AM"
$ws.Range("E43").Value = "FAIL"
$ws.Range("F43").Value = "1. Fast Return Time"
$ws.Range("G43").Value = "1. Static CSFB MO"
$ws.Range("H43").Value = "1. Reselection Time After CSFB Call Idle to LTE (ms): To achieve a faster return time for CSFB MO, please verify that the 2G serving cell belongs to the same site. Additionally, ensure that background download test files are running during the test to facilitate quicker reselection to LTE."

# Match the row height used by the similarly-wrapped row above it.
$ws.Rows.Item(43).RowHeight = 57

# Reflect the new active selection recorded in the saved workbook.
$ws.Range("F43").Select() | Out-Null
